# Scheduled market-data refresh: update cached Universalis price snapshots
# (columns H-N) across all 7 crafting-Leve sheets. Values are plain numeric
# literals pulled by the external runner; no formulas are involved.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 341.75
$ws.Range("I2").Value = 147.5
$ws.Range("K2").Value = 147.5
$ws.Range("M2").Value = -34.5
# Row 43
$ws.Range("H43").Value = 2899
$ws.Range("J43").Value = 2899
$ws.Range("L43").Value = 2899
$ws.Range("N43").Value = -3037
# Row 62
$ws.Range("H62").Value = 12115
$ws.Range("I62").Value = 8402.5
$ws.Range("K62").Value = 8402.5
$ws.Range("M62").Value = -7778.5
# Row 65
$ws.Range("H65").Value = 12115
$ws.Range("I65").Value = 8402.5
$ws.Range("K65").Value = 42012.5
$ws.Range("M65").Value = -38892.5
# Row 96
$ws.Range("H96").Value = 465.44446
$ws.Range("I96").Value = 362.7143
$ws.Range("J96").Value = 825
$ws.Range("K96").Value = 1088.1429
$ws.Range("L96").Value = 2475
$ws.Range("M96").Value = 284.8571000000002
$ws.Range("N96").Value = -5221
# Row 98
$ws.Range("H98").Value = 1173.1
$ws.Range("I98").Value = 1173.1
$ws.Range("K98").Value = 1173.1
$ws.Range("M98").Value = 324.9000000000001
# Row 116
$ws.Range("H116").Value = 4998.8335
$ws.Range("I116").Value = 4997.6665
$ws.Range("K116").Value = 4997.6665
$ws.Range("M116").Value = -1555.6665
# Row 122
$ws.Range("H122").Value = 1173.1
$ws.Range("I122").Value = 1173.1
$ws.Range("K122").Value = 3519.3
$ws.Range("M122").Value = -1069.3
# Row 131
$ws.Range("H131").Value = 1202
$ws.Range("I131").Value = 620.5
$ws.Range("K131").Value = 1861.5
$ws.Range("M131").Value = 3178.5
# Row 132
$ws.Range("H132").Value = 14238.474
$ws.Range("I132").Value = 16371.615
$ws.Range("K132").Value = 49114.845
$ws.Range("M132").Value = -46584.845
# Row 138
$ws.Range("H138").Value = 2142.923
$ws.Range("I138").Value = 1783.6666
$ws.Range("K138").Value = 5350.9998
$ws.Range("M138").Value = -210.9997999999996

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 3497.4
$ws.Range("J45").Value = 4199.8
$ws.Range("L45").Value = 4199.8
$ws.Range("N45").Value = -4953.8
# Row 61
$ws.Range("H61").Value = 6326.643
$ws.Range("I61").Value = 6162.727
$ws.Range("K61").Value = 6162.727
$ws.Range("M61").Value = -5950.727
# Row 122
$ws.Range("H122").Value = 1621.3334
$ws.Range("I122").Value = 1621.3334
$ws.Range("K122").Value = 4864.0002
$ws.Range("M122").Value = -2414.0002
# Row 132
$ws.Range("H132").Value = 2340.25
$ws.Range("I132").Value = 1216.25
$ws.Range("K132").Value = 3648.75
$ws.Range("M132").Value = -1118.75
# Row 136
$ws.Range("H136").Value = 6326.643
$ws.Range("I136").Value = 6162.727
$ws.Range("K136").Value = 18488.181
$ws.Range("M136").Value = -15938.181

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1774.5
$ws.Range("I16").Value = 1800
$ws.Range("J16").Value = 1749
$ws.Range("K16").Value = 1800
$ws.Range("L16").Value = 1749
$ws.Range("M16").Value = -1513
$ws.Range("N16").Value = -2323
# Row 58
$ws.Range("H58").Value = 3428.625
$ws.Range("J58").Value = 6313.3335
$ws.Range("L58").Value = 6313.3335
$ws.Range("N58").Value = -6719.3335
# Row 105
$ws.Range("H105").Value = 1400.7778
$ws.Range("I105").Value = 1572.75
$ws.Range("J105").Value = 25
$ws.Range("K105").Value = 1572.75
$ws.Range("L105").Value = 25
$ws.Range("M105").Value = 174.25
$ws.Range("N105").Value = -3519
# Row 107
$ws.Range("H107").Value = 232.2258
$ws.Range("I107").Value = 405
$ws.Range("J107").Value = 161.54546
$ws.Range("K107").Value = 405
$ws.Range("L107").Value = 161.54546
$ws.Range("M107").Value = 1515
$ws.Range("N107").Value = -4001.54546
# Row 113
$ws.Range("H113").Value = 1774.5
$ws.Range("I113").Value = 1800
$ws.Range("J113").Value = 1749
$ws.Range("K113").Value = 1800
$ws.Range("L113").Value = 1749
$ws.Range("M113").Value = 370
$ws.Range("N113").Value = -6089
# Row 122
$ws.Range("H122").Value = 1465.5264
$ws.Range("I122").Value = 1530.8462
$ws.Range("J122").Value = 1324
$ws.Range("K122").Value = 4592.5386
$ws.Range("L122").Value = 3972
$ws.Range("M122").Value = -2142.5386
$ws.Range("N122").Value = -8872
# Row 132
$ws.Range("H132").Value = 3842.6667
$ws.Range("I132").Value = 3405.6428
$ws.Range("K132").Value = 10216.9284
$ws.Range("M132").Value = -7686.928400000001
# Row 136
$ws.Range("H136").Value = 3428.625
$ws.Range("J136").Value = 6313.3335
$ws.Range("L136").Value = 18940.0005
$ws.Range("N136").Value = -24040.0005

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 8
$ws.Range("H8").Value = 357.41666
$ws.Range("I8").Value = 357.41666
$ws.Range("K8").Value = 1072.24998
$ws.Range("M8").Value = -933.2499800000001
# Row 32
$ws.Range("H32").Value = 1339.8
$ws.Range("I32").Value = 1424.75
$ws.Range("J32").Value = 1000
$ws.Range("K32").Value = 4274.25
$ws.Range("L32").Value = 3000
$ws.Range("M32").Value = -3991.25
$ws.Range("N32").Value = -3566
# Row 40
$ws.Range("H40").Value = 57.916668
$ws.Range("I40").Value = 40.2
$ws.Range("J40").Value = 146.5
$ws.Range("K40").Value = 160.8
$ws.Range("L40").Value = 586
$ws.Range("M40").Value = -91.80000000000001
$ws.Range("N40").Value = -724
# Row 107
$ws.Range("H107").Value = 774.25
$ws.Range("I107").Value = 245
$ws.Range("K107").Value = 735
$ws.Range("M107").Value = 1185

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 38
$ws.Range("H38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
# Row 46
$ws.Range("H46").Value = 8974
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
# Row 126
$ws.Range("H126").Value = 4588.1665
$ws.Range("I126").Value = 4201.4287
$ws.Range("K126").Value = 12604.2861
$ws.Range("M126").Value = -10134.2861
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 5763.769
$ws.Range("I7").Value = 5064.778
$ws.Range("J7").Value = 7336.5
$ws.Range("K7").Value = 5064.778
$ws.Range("L7").Value = 7336.5
$ws.Range("M7").Value = -4952.778
$ws.Range("N7").Value = -7560.5
# Row 22
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
# Row 27
$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
# Row 46
$ws.Range("H46").Value = 3648.72
$ws.Range("I46").Value = 1727.1
$ws.Range("J46").Value = 4929.8
$ws.Range("K46").Value = 1727.1
$ws.Range("L46").Value = 4929.8
$ws.Range("M46").Value = -1539.1
$ws.Range("N46").Value = -5305.8
# Row 104
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988
# Row 122
$ws.Range("H122").Value = 4054.3333
$ws.Range("I122").Value = 3513.3635
$ws.Range("K122").Value = 10540.0905
$ws.Range("M122").Value = -8090.0905
# Row 126
$ws.Range("H126").Value = 5763.769
$ws.Range("I126").Value = 5064.778
$ws.Range("J126").Value = 7336.5
$ws.Range("K126").Value = 15194.334
$ws.Range("L126").Value = 22009.5
$ws.Range("M126").Value = -12724.334
$ws.Range("N126").Value = -26949.5
# Row 132
$ws.Range("H132").Value = 3486.75
$ws.Range("I132").Value = 2688.8
$ws.Range("K132").Value = 8066.400000000001
$ws.Range("M132").Value = -5536.400000000001
# Row 136
$ws.Range("H136").Value = 3698.4546
$ws.Range("I136").Value = 3191.5557
$ws.Range("K136").Value = 9574.667099999999
$ws.Range("M136").Value = -7024.667099999999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 14
$ws.Range("H14").Value = 25000000
$ws.Range("I14").Value = 25000000
$ws.Range("K14").Value = 25000000
$ws.Range("M14").Value = -24999832
# Row 69
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 72
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 74
$ws.Range("H74").Value = 22077
$ws.Range("I74").Value = 19097.5
$ws.Range("K74").Value = 19097.5
$ws.Range("M74").Value = -18161.5
# Row 77
$ws.Range("H77").Value = 22077
$ws.Range("I77").Value = 19097.5
$ws.Range("K77").Value = 57292.5
$ws.Range("M77").Value = -52612.5
# Row 104
$ws.Range("H104").Value = 31831.4
$ws.Range("J104").Value = 31831.4
$ws.Range("L104").Value = 31831.4
$ws.Range("N104").Value = -38819.4
# Row 112
$ws.Range("H112").Value = 23591.334
$ws.Range("J112").Value = 23591.334
$ws.Range("L112").Value = 23591.334
$ws.Range("N112").Value = -26545.334
# Row 122
$ws.Range("H122").Value = 1303.7
$ws.Range("I122").Value = 1303.7
$ws.Range("K122").Value = 3911.1
$ws.Range("M122").Value = -1461.1
# Row 126
$ws.Range("H126").Value = 3461.0527
$ws.Range("I126").Value = 1404.2307
$ws.Range("K126").Value = 4212.6921
$ws.Range("M126").Value = -1742.6921
# Row 132
$ws.Range("H132").Value = 3076.4167
$ws.Range("I132").Value = 3076.4167
$ws.Range("K132").Value = 9229.250100000001
$ws.Range("M132").Value = -6699.250100000001
